$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1167).EntireRow.Insert()
$ws.Rows.Item(1168).EntireRow.Insert()

# Row 1167: Primera
$ws.Cells.Item(1167,1).Value = 3
$ws.Cells.Item(1167,2).Value = "Femacal de La Calera"
$ws.Cells.Item(1167,3).Value = "Coquimbo"
$ws.Cells.Item(1167,4).Value = 45223
$ws.Cells.Item(1167,5).Value = 5
$ws.Cells.Item(1167,6).Value = 100114014
$ws.Cells.Item(1167,7).Value = "Betarraga"
$ws.Cells.Item(1167,8).Value = "Sin especificar"
$ws.Cells.Item(1167,9).Value = "Primera"
$ws.Cells.Item(1167,10).Value = 1600
$ws.Cells.Item(1167,11).Value = 500
$ws.Cells.Item(1167,12).Value = 500
$ws.Cells.Item(1167,13).Value = 500
$ws.Cells.Item(1167,14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(1167,15).Value = "Provincia de Quillota"
$ws.Cells.Item(1167,16).Value = 125
$ws.Cells.Item(1167,17).Value = 4
$ws.Cells.Item(1167,18).Value = "Hortaliza"

# Row 1168: Segunda
$ws.Cells.Item(1168,1).Value = 3
$ws.Cells.Item(1168,2).Value = "Femacal de La Calera"
$ws.Cells.Item(1168,3).Value = "Coquimbo"
$ws.Cells.Item(1168,4).Value = 45223
$ws.Cells.Item(1168,5).Value = 5
$ws.Cells.Item(1168,6).Value = 100114014
$ws.Cells.Item(1168,7).Value = "Betarraga"
$ws.Cells.Item(1168,8).Value = "Sin especificar"
$ws.Cells.Item(1168,9).Value = "Segunda"
$ws.Cells.Item(1168,10).Value = 1600
$ws.Cells.Item(1168,11).Value = 400
$ws.Cells.Item(1168,12).Value = 400
$ws.Cells.Item(1168,13).Value = 400
$ws.Cells.Item(1168,14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(1168,15).Value = "Provincia de Quillota"
$ws.Cells.Item(1168,16).Value = 100
$ws.Cells.Item(1168,17).Value = 4
$ws.Cells.Item(1168,18).Value = "Hortaliza"

Write-Host "done"
